$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.35%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.35'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '5.81%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.133'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.03%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07645'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.60%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.268'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.42%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.615'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.13%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9073'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.01%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1135'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '10.49%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1797'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.60%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09115'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.44%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04232'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-4.33%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1050'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.46%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001259'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '2.40%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005729'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.24%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.45%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.65%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.751'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.35%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.33%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.54%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04061'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.01%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001267'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5.48%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004042'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.45%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001270'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.10%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003746'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02422'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-1.30%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05252'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.39%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007808'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.53%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1302'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.31%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006533'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-8.42%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001950'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.25%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007571'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.04%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.46%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006788'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '5.02%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.24%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06840'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1,424.77%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '40.27%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.24%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.24%'
